$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("AB12").Value = 160
$ws.Range("I12").Value = 14
$ws.Range("AB13").Value = 181
$ws.Range("I13").Value = 14
$ws.Range("AB14").Value = 206
$ws.Range("I14").Value = 18
$ws.Range("AB15").Value = 220
$ws.Range("I15").Value = 18
$ws.Range("AB16").Value = 236
$ws.Range("I16").Value = 18
$ws.Range("AB17").Value = 258
$ws.Range("I17").Value = 20
$ws.Range("AB18").Value = 284
$ws.Range("I18").Value = 20
$ws.Range("AB19").Value = 320
$ws.Range("I19").Value = 21
$ws.Range("AB20").Value = 331
$ws.Range("I20").Value = 21
$ws.Range("AB21").Value = 360
$ws.Range("I21").Value = 23
$ws.Range("AB22").Value = 384
$ws.Range("I22").Value = 27
$ws.Range("AB23").Value = 399
$ws.Range("I23").Value = 27
$ws.Range("AB24").Value = 433
$ws.Range("I24").Value = 28
$ws.Range("AB25").Value = 463
$ws.Range("I25").Value = 31
$ws.Range("AB26").Value = 520
$ws.Range("I26").Value = 33
$ws.Range("AB27").Value = 563
$ws.Range("I27").Value = 35
$ws.Range("AB28").Value = 623
$ws.Range("I28").Value = 38
$ws.Range("AB29").Value = 653
$ws.Range("I29").Value = 38
$ws.Range("AB30").Value = 687
$ws.Range("I30").Value = 39
$ws.Range("AB31").Value = 756
$ws.Range("I31").Value = 39
$ws.Range("AB32").Value = 817
$ws.Range("AB33").Value = 867
$ws.Range("J33").Value = 1

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("AB12").Value = 3
$ws.Range("I12").Value = 0
$ws.Range("AB13").Value = 3
$ws.Range("I13").Value = 0
$ws.Range("AB14").Value = 5
$ws.Range("I14").Value = 0
$ws.Range("AB15").Value = 5
$ws.Range("I15").Value = 0
$ws.Range("AB16").Value = 5
$ws.Range("I16").Value = 0
$ws.Range("AB17").Value = 5
$ws.Range("I17").Value = 0
$ws.Range("AB18").Value = 6
$ws.Range("I18").Value = 0
$ws.Range("AB19").Value = 8
$ws.Range("I19").Value = 0
$ws.Range("AB20").Value = 8
$ws.Range("I20").Value = 0
$ws.Range("AB21").Value = 8
$ws.Range("I21").Value = 0
$ws.Range("AB22").Value = 8
$ws.Range("I22").Value = 0
$ws.Range("AB23").Value = 8
$ws.Range("I23").Value = 0
$ws.Range("AB24").Value = 9
$ws.Range("I24").Value = 0
$ws.Range("AB25").Value = 10
$ws.Range("I25").Value = 0
$ws.Range("AB26").Value = 11
$ws.Range("I26").Value = 0
$ws.Range("AB27").Value = 11
$ws.Range("I27").Value = 0
$ws.Range("AB28").Value = 11
$ws.Range("I28").Value = 0
$ws.Range("AB29").Value = 11
$ws.Range("I29").Value = 0
$ws.Range("AB30").Value = 11
$ws.Range("I30").Value = 0
$ws.Range("AB31").Value = 12
$ws.Range("AB32").Value = 12
$ws.Range("AB33").Value = 12
$ws.Range("J33").Value = 1

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("AB4").Value = -5
$ws.Range("I4").Value = 3
$ws.Range("AB5").Value = -27
$ws.Range("I5").Value = 3
$ws.Range("AB6").Value = -27
$ws.Range("I6").Value = 5
$ws.Range("AB7").Value = -30
$ws.Range("I7").Value = 3
$ws.Range("AB8").Value = -39
$ws.Range("AB9").Value = -41
$ws.Range("AB10").Value = -55
$ws.Range("I10").Value = 5
$ws.Range("AB11").Value = -53
$ws.Range("I11").Value = 8
$ws.Range("AB12").Value = -56
$ws.Range("I12").Value = 7
$ws.Range("AB13").Value = -65
$ws.Range("I13").Value = 2
$ws.Range("AB14").Value = -59
$ws.Range("I14").Value = 5
$ws.Range("AB15").Value = -63
$ws.Range("AB16").Value = -63
$ws.Range("AB17").Value = -74
$ws.Range("I17").Value = 1
$ws.Range("AB18").Value = -78
$ws.Range("I18").Value = 1
$ws.Range("AB19").Value = -73
$ws.Range("I19").Value = 1
$ws.Range("AB20").Value = -73
$ws.Range("I20").Value = -1
$ws.Range("AB21").Value = -78
$ws.Range("I21").Value = -1
$ws.Range("AB22").Value = -78
$ws.Range("AB23").Value = -76
$ws.Range("AB24").Value = -71
$ws.Range("I24").Value = 4
$ws.Range("AB25").Value = -76
$ws.Range("I25").Value = 4
$ws.Range("AB26").Value = -75
$ws.Range("AB27").Value = -80
$ws.Range("I27").Value = 3
$ws.Range("AB28").Value = -77
$ws.Range("I28").Value = 2
$ws.Range("AB29").Value = -77
$ws.Range("AB30").Value = -78
$ws.Range("AB31").Value = -69
$ws.Range("I31").Value = 2
$ws.Range("AB32").Value = -68
$ws.Range("AB33").Value = -67
$ws.Range("J33").Value = 0

$ws = $wb.Worksheets.Item("ICU")
$ws.Range("I19").Value = 1
$ws.Range("AB20").Value = -21
$ws.Range("I20").Value = 2
$ws.Range("AB21").Value = -19
$ws.Range("I21").Value = 3
$ws.Range("AB22").Value = -20
$ws.Range("AB23").Value = -21
$ws.Range("AB24").Value = -21
$ws.Range("I24").Value = 3
$ws.Range("AB25").Value = -20
$ws.Range("I25").Value = 4
$ws.Range("AB26").Value = -21
$ws.Range("AB27").Value = -20
$ws.Range("I27").Value = 5
$ws.Range("AB28").Value = -21
$ws.Range("I28").Value = 4
$ws.Range("AB29").Value = -21
$ws.Range("AB30").Value = -21
$ws.Range("AB31").Value = -21
$ws.Range("I31").Value = 3
$ws.Range("AB32").Value = -21
$ws.Range("AB33").Value = -21

$ws = $wb.Worksheets.Item("Ventilated")
$ws.Range("I19").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("AB24").Value = -5
$ws.Range("I24").Value = 1
$ws.Range("AB25").Value = -4
$ws.Range("I25").Value = 2
$ws.Range("AB26").Value = -4
$ws.Range("AB27").Value = -3
$ws.Range("I27").Value = 2
$ws.Range("AB28").Value = -4
$ws.Range("I28").Value = 1
$ws.Range("AB29").Value = -5
$ws.Range("AB30").Value = -5
$ws.Range("AB31").Value = -4
$ws.Range("I31").Value = 1
$ws.Range("AB32").Value = -5
$ws.Range("AB33").Value = -5
